$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header: "Skills Match" -> "Skills/Interest Match"
$ws.Range("A12").Value = "Skills/Interest Match"

# "Career Prospects Match" section (rows 28-36): fill in the rating columns
# B-F for every indicator row, rename "Accounting, Auditing" to
# "Accounting, Financial Auditing" (row 30), and append a new row 36
# ("Swiss Organization").
$rows = @(
    "28|Finance, Banking|high|medium-high|medium|medium|low-medium",
    "29|Marketing, Communication|medium-high|medium-high|low-medium|low-medium|medium-high",
    "30|Accounting, Financial Auditing|high|medium-high|low|low|medium",
    "31|Business Analytics, Business Intelligence|low-medium|low-medium|high|high|low",
    "32|Software Development|low|low|high|high|low",
    "33|HRM, Leadership|medium-high|medium-high|low|low|medium-high",
    "34|Entrepreneurship, Starting own Business|medium-high|medium-high|medium-high|medium-high|medium-high",
    "35|International Organizations|medium|high|medium|high|high",
    "36|Swiss Organization|high|medium|high|medium|high"
)

foreach ($line in $rows) {
    $p = $line.Split("|")
    $r = [int]$p[0]
    $ws.Cells.Item($r, 1).Value = $p[1]
    $ws.Cells.Item($r, 2).Value = $p[2]
    $ws.Cells.Item($r, 3).Value = $p[3]
    $ws.Cells.Item($r, 4).Value = $p[4]
    $ws.Cells.Item($r, 5).Value = $p[5]
    $ws.Cells.Item($r, 6).Value = $p[6]
}

# Freeze the header row and leave the selection on the next free row
# beneath the newly added data (matches the saved view state).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A37").Select()
